# Update cryptocurrency price values (column D) on the active worksheet.
# Values are stored as text strings in the workbook, so we set them
# explicitly as strings to preserve exact formatting (e.g. trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "246.36"
    "D3"  = "21.99"
    "D4"  = "5.414"
    "D6"  = "3.392"
    "D7"  = "6.339"
    "D8"  = "0.8091"
    "D9"  = "0.9519"
    "D10" = "0.1423"
    "D11" = "0.07505"
    "D12" = "0.03187"
    "D13" = "0.03015"
    "D14" = "4.138"
    "D15" = "0.09408"
    "D16" = "0.001594"
    "D17" = "0.04823"
    "D18" = "0.0005893"
    "D19" = "0.006191"
    "D20" = "0.004112"
    "D21" = "0.0009971"
    "D23" = "3.767"
    "D24" = "2.230"
    "D27" = "0.0001291"
    "D40" = "0.03887"
    "D41" = "0.006351"
    "D42" = "0.1073"
    "D43" = "0.003000"
    "D44" = "0.006693"
    "D45" = "0.00005596"
    "D47" = "0.3802"
    "D48" = "0.1477"
    "D49" = "0.00002100"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
